$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 237 — every existing row from 237..255 shifts
# down by one (becoming 238..256), matching the dimension growing from
# A1:R255 to A1:R256.
$ws.Rows.Item(237).Insert()

# Populate the newly inserted row 237 with the new weekly price record.
$ws.Cells.Item(237, 1).Value = 10
$ws.Cells.Item(237, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(237, 3).Value = "La Araucanía"
$ws.Cells.Item(237, 4).Value = 44578
$ws.Cells.Item(237, 5).Value = 9
$ws.Cells.Item(237, 6).Value = 100112009
$ws.Cells.Item(237, 7).Value = "Acelga"
$ws.Cells.Item(237, 8).Value = "Sin especificar"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 50
$ws.Cells.Item(237, 11).Value = 8000
$ws.Cells.Item(237, 12).Value = 8000
$ws.Cells.Item(237, 13).Value = 8000
$ws.Cells.Item(237, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(237, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(237, 16).Value = 667
$ws.Cells.Item(237, 17).Value = 12
$ws.Cells.Item(237, 18).Value = "Hortaliza"
